$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.04042238960271747
$ws.Range("J2").Value = 0.04042238960271747
$ws.Range("M2").Value = 0.096887
$ws.Range("N2").Value = 0.290661
$ws.Range("O2").Value = 0.009776580706310958
$ws.Range("P2").Value = 0.009776580706310958
$ws.Range("Q2").Value = 0.002033044512333333
$ws.Range("R2").Value = 0.018297400611
$ws.Range("S2").Value = 0.0003951927542929122
$ws.Range("T2").Value = 0.0003951927542929122

# Row 3
$ws.Range("I3").Value = 0.04042238960271747
$ws.Range("J3").Value = 0.04042238960271747
$ws.Range("O3").Value = 0.4623470245018782
$ws.Range("P3").Value = 0.4623470245018782
$ws.Range("S3").Value = 0.01868917155607208
$ws.Range("T3").Value = 0.01868917155607208

# Row 4
$ws.Range("I4").Value = 0.04042238960271747
$ws.Range("J4").Value = 0.04042238960271747
$ws.Range("M4").Value = 5.076459666666667
$ws.Range("N4").Value = 15.229379
$ws.Range("O4").Value = 0.51225053550527
$ws.Range("P4").Value = 0.51225053550527
$ws.Range("Q4").Value = 0.1065227374921111
$ws.Range("R4").Value = 0.9587046374289999
$ws.Range("S4").Value = 0.02070639072039468
$ws.Range("T4").Value = 0.02070639072039468

# Row 5
$ws.Range("I5").Value = 0.04042238960271747
$ws.Range("J5").Value = 0.04042238960271747
$ws.Range("M5").Value = 0.154854
$ws.Range("N5").Value = 0.464562
$ws.Range("O5").Value = 0.01562585928654078
$ws.Range("P5").Value = 0.01562585928654078
$ws.Range("Q5").Value = 0.003249404717999999
$ws.Range("R5").Value = 0.029244642462
$ws.Range("S5").Value = 0.0006316345719577924
$ws.Range("T5").Value = 0.0006316345719577924

# Row 6
$ws.Range("G6").Value = 0.4981263333333333
$ws.Range("H6").Value = 1.494379
$ws.Range("I6").Value = 0.9595776103972825
$ws.Range("J6").Value = 0.9595776103972825
$ws.Range("M6").Value = 0.096887
$ws.Range("N6").Value = 0.290661
$ws.Range("O6").Value = 0.009776580706310958
$ws.Range("P6").Value = 0.009776580706310958
$ws.Range("Q6").Value = 0.04826196605766666
$ws.Range("R6").Value = 0.434357694519
$ws.Range("S6").Value = 0.009381387952018046
$ws.Range("T6").Value = 0.009381387952018046

# Row 7
$ws.Range("G7").Value = 0.4981263333333333
$ws.Range("H7").Value = 1.494379
$ws.Range("I7").Value = 0.9595776103972825
$ws.Range("J7").Value = 0.9595776103972825
$ws.Range("O7").Value = 0.4623470245018782
$ws.Range("P7").Value = 0.4623470245018782
$ws.Range("Q7").Value = 2.282370194005444
$ws.Range("R7").Value = 20.541331746049
$ws.Range("S7").Value = 0.4436578529458061
$ws.Range("T7").Value = 0.4436578529458061

# Row 8
$ws.Range("G8").Value = 0.4981263333333333
$ws.Range("H8").Value = 1.494379
$ws.Range("I8").Value = 0.9595776103972825
$ws.Range("J8").Value = 0.9595776103972825
$ws.Range("M8").Value = 5.076459666666667
$ws.Range("N8").Value = 15.229379
$ws.Range("O8").Value = 0.51225053550527
$ws.Range("P8").Value = 0.51225053550527
$ws.Range("Q8").Value = 2.528718240071222
$ws.Range("R8").Value = 22.758464160641
$ws.Range("S8").Value = 0.4915441447848753
$ws.Range("T8").Value = 0.4915441447848753

# Row 9
$ws.Range("G9").Value = 0.4981263333333333
$ws.Range("H9").Value = 1.494379
$ws.Range("I9").Value = 0.9595776103972825
$ws.Range("J9").Value = 0.9595776103972825
$ws.Range("M9").Value = 0.154854
$ws.Range("N9").Value = 0.464562
$ws.Range("O9").Value = 0.01562585928654078
$ws.Range("P9").Value = 0.01562585928654078
$ws.Range("Q9").Value = 0.07713685522199999
$ws.Range("R9").Value = 0.6942316969979999
$ws.Range("S9").Value = 0.01499422471458299
$ws.Range("T9").Value = 0.01499422471458299
